# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "CasesTab" query (row 2, column B) used to contain an extra trailing
# line returning a `Cohort` column. That line is removed from the query
# text. Because the modified text is a brand-new unique value, Excel will
# allocate it a new shared-string slot; the now-unreferenced previous text
# is dropped from the shared string table, which in turn shifts the
# "FilesTab" query (row 4, column B) into the slot vacated by the old
# "CasesTab" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n`tWHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in ['T2N1M0', 'T3N0M0', 'T3N0M1', 'T3N1M0'] OPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value2 = $newCasesQuery

# The row height of row 2 was tied to the old (longer) text wrapping to
# more lines; with the shorter text it now matches the height already
# used by the other query rows (290).
$ws.Rows.Item(2).RowHeight = 290

# Update the window view to match the new focal point of the edit: the
# sheet had scrolled to show row 4 / cell C4 selected; now it shows the
# top of the sheet (B2) scrolled into the top-left corner, with B2 selected.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B2").Select() | Out-Null
